# New crime data collected: refresh the weekly CompStat sheet (26th Precinct)
# for the week of 7/3/2023 - 7/9/2023, update the Police Commissioner byline,
# and roll the report numbers forward (Volume 30, Number 27) along with all
# of the underlying weekly / 28-day / year-to-date / 2-year crime statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header text updates -------------------------------------------------
# Police Commissioner byline
$ws.Range("M6").Value = "Edward A. Caban"
# "Volume 30   Number  26" -> "...27"
$ws.Range("A8").Value = "Volume 30   Number  27"
# "Report Covering the Week  6/26/2023  Through  7/2/2023" -> next week
$ws.Range("C9").Value = "Report Covering the Week  7/3/2023  Through  7/9/2023"

# --- simple numeric value updates (style unchanged) -----------------------
$numericUpdates = @{
    "C16" = 4
    "D16" = 4
    "E16" = 0
    "F16" = 15
    "G16" = 9
    "H16" = 66.666666666666
    "I16" = 70
    "J16" = 51
    "K16" = 37.254901960784
    "L16" = 100
    "M16" = -14.634146341463
    "N16" = -83.009708737864
    "C17" = 2
    "D17" = 6
    "E17" = -66.666666666666
    "F17" = 6
    "G17" = 16
    "H17" = -62.5
    "I17" = 85
    "J17" = 75
    "K17" = 13.333333333333
    "L17" = 16.438356164383
    "M17" = 32.8125
    "N17" = -44.805194805194
    "E18" = 200
    "F18" = 5
    "H18" = 0
    "I18" = 49
    "J18" = 64
    "K18" = -23.4375
    "L18" = 68.965517241379
    "M18" = 22.5
    "N18" = -85.106382978723
    "C19" = 7
    "D19" = 6
    "E19" = 16.666666666666
    "F19" = 33
    "G19" = 24
    "H19" = 37.5
    "I19" = 210
    "J19" = 176
    "K19" = 19.318181818181
    "L19" = 87.5
    "M19" = 51.079136690647
    "N19" = -38.775510204081
    "G20" = 5
    "H20" = 100
    "J20" = 39
    "K20" = 66.666666666666
    "L20" = 282.352941176471
    "M20" = 550
    "N20" = -69.194312796208
    "C21" = 16
    "D21" = 20
    "E21" = -20
    "F21" = 69
    "G21" = 59
    "H21" = 16.949152542372
    "I21" = 483
    "J21" = 413
    "K21" = 16.949152542372
    "L21" = 78.888888888888
    "M21" = 40.406976744186
    "N21" = -67.165193745751
    "C22" = 1
    "F22" = 4
    "I22" = 21
    "K22" = 110
    "L22" = 40
    "M22" = 50
    "C23" = 2
    "D23" = 4
    "F23" = 10
    "H23" = 11.111111111111
    "I23" = 76
    "J23" = 58
    "K23" = 31.034482758620
    "L23" = 38.181818181818
    "M23" = 61.702127659574
    "C24" = 10
    "D24" = 8
    "E24" = 25
    "F24" = 46
    "G24" = 31
    "H24" = 48.387096774193
    "I24" = 276
    "J24" = 263
    "K24" = 4.942965779467
    "L24" = 35.960591133004
    "M24" = -4.827586206896
    "C25" = 3
    "D25" = 2
    "E25" = 50
    "F25" = 13
    "G25" = 15
    "H25" = -13.333333333333
    "I25" = 127
    "J25" = 127
    "K25" = 0
    "L25" = 39.560439560439
    "M25" = -23.030303030303
    "G27" = 4
    "H27" = 25
    "L27" = 71.428571428571
    "N28" = -60
    "N29" = -57.142857142857
}
foreach ($addr in $numericUpdates.Keys) {
    $ws.Range($addr).Value = $numericUpdates[$addr]
}

# --- cells converting FROM a numeric style TO the shared-string
#     placeholder style (s="14", used for "0" / "***.*" blanks) -----------
# Copying a cell that already carries the placeholder keeps the same style
# index and shared-string reference instead of minting a new style/string.

# '0' placeholder (shared string index 20)
$zeroPlaceholderCells = @("G15", "C20", "G26", "C27", "D28", "D29", "G30")
foreach ($addr in $zeroPlaceholderCells) {
    $ws.Range("C14").Copy($ws.Range($addr))
}

# '***.*' placeholder (shared string index 21)
$naPlaceholderCells = @("H15", "H26", "E28", "E29", "H30")
foreach ($addr in $naPlaceholderCells) {
    $ws.Range("E14").Copy($ws.Range($addr))
}

# --- cells converting FROM the shared-string placeholder style (14) TO a
#     numeric style, now that they carry a real figure again --------------
$toStyle15 = @{
    "C18" = 3
    "D20" = 3
}
foreach ($addr in $toStyle15.Keys) {
    $ws.Range("J14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $toStyle15[$addr]
}

$toStyle16 = @{
    "E20" = -100
}
foreach ($addr in $toStyle16.Keys) {
    $ws.Range("K14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $toStyle16[$addr]
}
